$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.808.04"
$ws.Range("E2").Value = "  +4.31%  "

$ws.Range("D3").Value = "2.254.22"
$ws.Range("E3").Value = "  +3.55%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.79"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.672"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +18.23%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +10.11%  "

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.24"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.92%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0961"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +8.59%  "

$ws.Range("E14").Value = "  +1.44%  "

$ws.Range("D15").Value = "2.584.29"
$ws.Range("E15").Value = "  +3.50%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.77"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.03%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.881"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.86%  "

$ws.Range("D18").Value = "2.258.86"
$ws.Range("E18").Value = "  +4.15%  "

$ws.Range("D19").Value = "42.738.75"
$ws.Range("E19").Value = "  +4.39%  "

$ws.Range("D20").Value = "0.0₃0990"
$ws.Range("E20").Value = "  +5.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.30"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.93"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.33%  "

$ws.Range("E24").Value = "  +0.81%  "

$ws.Range("E25").Value = "  +5.37%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.61"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.41%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").Value = "  +0.59%  "

$ws.Range("E29").Value = "  -1.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.43"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.99"
$ws.Range("D32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.50"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +15.00%  "

$ws.Range("E34").Value = "  +7.16%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.22"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +22.07%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0787"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.127"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.38"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.72"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.10%  "

$ws.Range("E40").Value = "  +6.94%  "

$ws.Range("E41").Value = "  +6.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.75"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +11.50%  "

$ws.Range("E43").Value = "  +6.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.81"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.00%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.09"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.51%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.202"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.83"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("E48").Value = "  +4.37%  "

$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.17"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("E51").Value = "  +4.04%  "
